$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared-string rich text cells) -------------------
# A8: "Volume 30   Number  33" -> "Volume 30   Number  34"
$ws.Range("A8").Value = "Volume 30   Number  34"

# C9: "Report Covering the Week  8/14/2023  Through  8/20/2023"
#  -> "Report Covering the Week  8/21/2023  Through  8/27/2023"
$ws.Range("C9").Value = "Report Covering the Week  8/21/2023  Through  8/27/2023"

# --- Row 16 (Robbery) -------------------------------------------------------
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("F16").Value = 6
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 44
$ws.Range("J16").Value = 59
$ws.Range("K16").Value = -25.423728813559
$ws.Range("L16").Value = 193.333333333333
$ws.Range("M16").Value = -18.518518518518
$ws.Range("N16").Value = -79.816513761467

# --- Row 17 (Fel. Assault) --------------------------------------------------
$ws.Range("C17").Value = 2
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 12.5
$ws.Range("I17").Value = 64
$ws.Range("J17").Value = 58
$ws.Range("K17").Value = 10.344827586206
$ws.Range("L17").Value = 39.130434782608
$ws.Range("M17").Value = 88.235294117647
$ws.Range("N17").Value = -12.328767123287

# --- Row 18 (Burglary) ------------------------------------------------------
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -57.142857142857
$ws.Range("F18").Value = 22
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = 10
$ws.Range("I18").Value = 203
$ws.Range("J18").Value = 170
$ws.Range("K18").Value = 19.411764705882
$ws.Range("L18").Value = 40.972222222222
$ws.Range("M18").Value = 26.086956521739
$ws.Range("N18").Value = -69.701492537313

# --- Row 19 (Gr. Larceny) ---------------------------------------------------
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 22.222222222222
$ws.Range("F19").Value = 50
$ws.Range("G19").Value = 56
$ws.Range("H19").Value = -10.714285714285
$ws.Range("I19").Value = 431
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = 7.75
$ws.Range("L19").Value = 78.838174273858
$ws.Range("M19").Value = 84.188034188034
$ws.Range("N19").Value = 24.566473988439

# --- Row 20 (G.L.A.) ---------------------------------------------------------
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 109
$ws.Range("J20").Value = 67
$ws.Range("K20").Value = 62.686567164179
$ws.Range("L20").Value = 136.95652173913
$ws.Range("M20").Value = 15.957446808510
$ws.Range("N20").Value = -94.772182254196

# --- Row 21 (TOTAL) ---------------------------------------------------------
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = -4.761904761904
$ws.Range("F21").Value = 106
$ws.Range("G21").Value = 101
$ws.Range("H21").Value = 4.950495049504
$ws.Range("I21").Value = 861
$ws.Range("J21").Value = 756
$ws.Range("K21").Value = 13.888888888888
$ws.Range("L21").Value = 73.939393939393
$ws.Range("M21").Value = 48.192771084337
$ws.Range("N21").Value = -74.691358024691

# --- Row 24 (Petit Larceny) -------------------------------------------------
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = -6.666666666666
$ws.Range("F24").Value = 43
$ws.Range("G24").Value = 56
$ws.Range("H24").Value = -23.214285714285
$ws.Range("I24").Value = 376
$ws.Range("J24").Value = 506
$ws.Range("K24").Value = -25.691699604743
$ws.Range("L24").Value = 5.617977528089
$ws.Range("M24").Value = 22.077922077922

# --- Row 25 (Misd. Assault) -------------------------------------------------
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 16
$ws.Range("G25").Value = 15
$ws.Range("H25").Value = 6.666666666666
$ws.Range("I25").Value = 140
$ws.Range("J25").Value = 156
$ws.Range("K25").Value = -10.256410256410
$ws.Range("L25").Value = 55.555555555555
$ws.Range("M25").Value = 26.126126126126

# --- Row 27 (Other Sex Crimes): C27 2 -> blank "0" placeholder -------------
$ws.Range("C27").Value = "0"

"done"
